$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: the worksheet that used to be "总计" (sheetId 6) is repurposed to
# hold the new "2022-Q1" fund-holding detail data (same layout as the other
# quarterly sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值/仓位排名).
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# template sheet used only to copy the header/index-column formatting (style)
$tmpl = $wb.Worksheets.Item("2021-Q4")

function Set-TextValue($range, $value) {
    # Write $value as literal text (even if it looks like a number), without
    # leaving a lingering custom number-format style on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Header row (copy bold/border/centered style from the template, then set text)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # column B..H
    $tmpl.Cells.Item(1, $col).Copy()
    $q1.Cells.Item(1, $col).PasteSpecial(-4122)
    $q1.Cells.Item(1, $col).Value = $headers[$i]
}

$data = @(
    @("590008", "中邮战略新兴产业混合", "9.00", "88.88", "4.14", "0.3726", 5),
    @("012421", "华夏优加生活混合A", "8.67", "92.98", "2.77", "0.2402", 7),
    @("004965", "泓德致远混合A", "21.03", "46.32", "1.06", "0.2229", 10),
    @("004966", "泓德致远混合C", "2.92", "46.32", "1.06", "0.0310", 10),
    @("012422", "华夏优加生活混合C", "0.17", "92.98", "2.77", "0.0047", 7)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = 2 + $r
    $rowData = $data[$r]

    # A column: index, styled like the template's index column
    $tmpl.Cells.Item(2, 1).Copy()
    $q1.Cells.Item($row, 1).PasteSpecial(-4122)
    $q1.Cells.Item($row, 1).Value = $r

    Set-TextValue $q1.Cells.Item($row, 2) $rowData[0]
    Set-TextValue $q1.Cells.Item($row, 3) $rowData[1]
    Set-TextValue $q1.Cells.Item($row, 4) $rowData[2]
    Set-TextValue $q1.Cells.Item($row, 5) $rowData[3]
    Set-TextValue $q1.Cells.Item($row, 6) $rowData[4]
    Set-TextValue $q1.Cells.Item($row, 7) $rowData[5]

    $q1.Cells.Item($row, 8).Value = $rowData[6]
}

# ---------------------------------------------------------------------------
# Step 2: add a brand-new "总计" worksheet right after "2022-Q1" holding the
# updated totals-per-quarter table (with the new 2022-Q1 row inserted on top).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# match the page margins used throughout the rest of the workbook
$total.PageSetup.LeftMargin = 0.75 * 72
$total.PageSetup.RightMargin = 0.75 * 72
$total.PageSetup.TopMargin = 1 * 72
$total.PageSetup.BottomMargin = 1 * 72
$total.PageSetup.HeaderMargin = 0.5 * 72
$total.PageSetup.FooterMargin = 0.5 * 72

$tmpl.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalData = @(
    @("2022-Q1", 5, 0.87),
    @("2021-Q4", 11, 4.12),
    @("2021-Q3", 20, 5.68),
    @("2021-Q2", 10, 5.16),
    @("2021-Q1", 17, 5.87),
    @("2020-Q4", 17, 2.83)
)

for ($r = 0; $r -lt $totalData.Length; $r++) {
    $row = 2 + $r
    $rowData = $totalData[$r]

    $tmpl.Cells.Item(2, 1).Copy()
    $total.Cells.Item($row, 1).PasteSpecial(-4122)
    $total.Cells.Item($row, 1).Value = $r

    $total.Cells.Item($row, 2).Value = $rowData[0]
    $total.Cells.Item($row, 3).Value = $rowData[1]
    $total.Cells.Item($row, 4).Value = $rowData[2]
}
